$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 47
$ws1.Range("H47").Value = 205.2
$ws1.Range("I47").Value = 152.3
$ws1.Range("M47").Value = 969.41

# Row 55
$ws1.Range("E55").Value = 414.08
$ws1.Range("P55").Value = 550.63

# Row 59 - "N de 57" progress counters
$ws1.Range("E59").Value = "3 de 57"
$ws1.Range("H59").Value = "2 de 57"
$ws1.Range("I59").Value = "2 de 57"
$ws1.Range("M59").Value = "5 de 57"
$ws1.Range("P59").Value = "3 de 57"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F47").Value = 1202.67
$ws2.Range("F55").Value = 964.71
$ws2.Range("F59").Value = 10462.96

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Widen column E (5th column) from 22 to 23.
# Note: the COM ColumnWidth property is offset from the raw OOXML <col width>
# value by 5/6 (~0.8333333) character units on this engine, so subtract that
# offset to land exactly on a stored width of 23.
$ws3.Columns.Item(5).ColumnWidth = 22.166666666666668

# Row 4 - FREGADEROS DE COCINA
$ws3.Range("D4").Value = 704.77
$ws3.Range("E4").Value = 77.64716394895902
$ws3.Range("F4").Value = 0.9007598918752447

# Row 6 - INODOROS
$ws3.Range("D6").Value = 430.8
$ws3.Range("E6").Value = 1895.26694516821
$ws3.Range("F6").Value = 0.1852053316414101

# Row 7 - LAVABOS
$ws3.Range("D7").Value = 542.67
$ws3.Range("E7").Value = 344.0410162875741
$ws3.Range("F7").Value = 0.6120032231831478

# Row 8 - NO RESURTIBLES
$ws3.Range("D8").Value = 777.21
$ws3.Range("E8").Value = -328.40837082797
$ws3.Range("F8").Value = 1.731745050555705

# Row 11 - PIEDRA SINTERIZADA
$ws3.Range("D11").Value = 2231.76
$ws3.Range("E11").Value = 17341.3002492497
$ws3.Range("F11").Value = 0.1140220267847769

# Row 12 - PORCELANATO
$ws3.Range("D12").Value = 6196.67
$ws3.Range("E12").Value = 42427.39
$ws3.Range("F12").Value = 0.1274404070741933

# Row 14 - TOTAL
$ws3.Range("D14").Value = 14031.58
$ws3.Range("E14").Value = 85866.41284188785
$ws3.Range("F14").Value = 0.1404590783140987
